$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(1.0118548668990655, 0.7726613778350778, 1.1562842043897095, 2.938328189687746)"
$ws.Range("C2").Value = "NIG(1.3760349876318032, 1.0573586881370223, 4.318090560111765, 5.778994955972886)"
$ws.Range("D2").Value = "NCT(3.099554074038097, 1.419277526749148, -0.45980177100041486, 2.7171635131087912)"
$ws.Range("E2").Value = "JSU(-1.2973905136434918, 1.1963036930864677, 1.9029651265620724, 5.097354282656984)"
